$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.09594307528308157
$ws.Range("J2").Value = 0.09594307528308157
$ws.Range("M2").Value = 1.619868333333333
$ws.Range("N2").Value = 4.859605
$ws.Range("O2").Value = 0.1089327058120143
$ws.Range("P2").Value = 0.1089327058120143
$ws.Range("Q2").Value = 0.9751137604850001
$ws.Range("R2").Value = 8.776023844365001
$ws.Range("S2").Value = 0.01045133879451187
$ws.Range("T2").Value = 0.01045133879451187

# Row 3
$ws.Range("I3").Value = 0.09594307528308157
$ws.Range("J3").Value = 0.09594307528308157
$ws.Range("O3").Value = 0.1655705935257241
$ws.Range("P3").Value = 0.1655705935257241
$ws.Range("S3").Value = 0.01588535191930305
$ws.Range("T3").Value = 0.01588535191930305

# Row 4
$ws.Range("I4").Value = 0.09594307528308157
$ws.Range("J4").Value = 0.09594307528308157
$ws.Range("M4").Value = 7.682722666666667
$ws.Range("N4").Value = 23.048168
$ws.Range("O4").Value = 0.5166467859527435
$ws.Range("P4").Value = 0.5166467859527435
$ws.Range("Q4").Value = 4.624776246376
$ws.Range("R4").Value = 41.622986217384
$ws.Range("S4").Value = 0.0495686814794262
$ws.Range("T4").Value = 0.0495686814794262

# Row 5
$ws.Range("I5").Value = 0.09594307528308157
$ws.Range("J5").Value = 0.09594307528308157
$ws.Range("M5").Value = 3.105673
$ws.Range("N5").Value = 9.317019
$ws.Range("O5").Value = 0.2088499147095181
$ws.Range("P5").Value = 0.2088499147095181
$ws.Range("Q5").Value = 1.869525081483
$ws.Range("R5").Value = 16.825725733347
$ws.Range("S5").Value = 0.02003770308984046
$ws.Range("T5").Value = 0.02003770308984046

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.672280666666667
$ws.Range("H6").Value = 17.016842
$ws.Range("I6").Value = 0.9040569247169185
$ws.Range("J6").Value = 0.9040569247169185
$ws.Range("M6").Value = 1.619868333333333
$ws.Range("N6").Value = 4.859605
$ws.Range("O6").Value = 0.1089327058120143
$ws.Range("P6").Value = 0.1089327058120143
$ws.Range("Q6").Value = 9.188347829712223
$ws.Range("R6").Value = 82.69513046741001
$ws.Range("S6").Value = 0.09848136701750244
$ws.Range("T6").Value = 0.09848136701750243

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.672280666666667
$ws.Range("H7").Value = 17.016842
$ws.Range("I7").Value = 0.9040569247169185
$ws.Range("J7").Value = 0.9040569247169185
$ws.Range("O7").Value = 0.1655705935257241
$ws.Range("P7").Value = 0.1655705935257241
$ws.Range("Q7").Value = 13.965688195716
$ws.Range("R7").Value = 125.691193761444
$ws.Range("S7").Value = 0.1496852416064211
$ws.Range("T7").Value = 0.1496852416064211

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.672280666666667
$ws.Range("H8").Value = 17.016842
$ws.Range("I8").Value = 0.9040569247169185
$ws.Range("J8").Value = 0.9040569247169185
$ws.Range("M8").Value = 7.682722666666667
$ws.Range("N8").Value = 23.048168
$ws.Range("O8").Value = 0.5166467859527435
$ws.Range("P8").Value = 0.5166467859527435
$ws.Range("Q8").Value = 43.57855924949511
$ws.Range("R8").Value = 392.207033245456
$ws.Range("S8").Value = 0.4670781044733173
$ws.Range("T8").Value = 0.4670781044733173

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.672280666666667
$ws.Range("H9").Value = 17.016842
$ws.Range("I9").Value = 0.9040569247169185
$ws.Range("J9").Value = 0.9040569247169185
$ws.Range("M9").Value = 3.105673
$ws.Range("N9").Value = 9.317019
$ws.Range("O9").Value = 0.2088499147095181
$ws.Range("P9").Value = 0.2088499147095181
$ws.Range("Q9").Value = 17.61624891488866
$ws.Range("R9").Value = 158.546240233998
$ws.Range("S9").Value = 0.1888122116196776
$ws.Range("T9").Value = 0.1888122116196776
